$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All fitness values in column C (rows 2 through 252) are updated to 7293
$ws.Range("C2:C252").Value = 7293
